{"js": "// Remove the \"Altman, Sam. 2025. 'Reflections.' Blog post.\" bullet from the\n// \"Consult as desired:\" sub-list. (The other two bullets - \"Consciousness\"\n// and \"Artificial Intelligence\" - are left untouched.)\nconst results = context.document.body.search(\"Altman, Sam. 2025.\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the \"Altman, Sam. 2025.\" bullet to remove.');\n}\n\nconst hit = results.items[0];\nconst targetParagraph = hit.paragraphs.getFirst();\ntargetParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the \"Altman, Sam. 2025. 'Reflections.' Blog post.\" bullet from the\n# \"Consult as desired:\" sub-list. (The other two bullets - \"Consciousness\"\n# and \"Artificial Intelligence\" - are left untouched.)\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"Altman, Sam. 2025.\")\nif (-not $found) {\n    throw \"Could not find the 'Altman, Sam. 2025.' bullet to remove.\"\n}\n\n# Expand the found range to the whole paragraph (wdParagraph = 4) so the\n# paragraph mark is included, then delete it outright.\n$range.Expand(4)\n$range.Delete()\n"}
